# "add parameters to birt" - turn the filled-in demo/sample sheet back into a
# blank BIRT report template: wipe the sample header labels and the sample
# data-row formatting/placeholder row, leaving only the "[Table1]" marker
# cell that BIRT replaces with real tabular data at report-generation time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was the sample table header ("Город" / "e-mail" / "баланс" / "t").
# Drop the text but keep the bold header look, minus its border (the
# template no longer draws a boxed header).
$ws.Range("B2:E2").ClearContents()
$ws.Range("B2:E2").Borders.LineStyle = 0

# Row 3 held one sample data row plus the "[Table1]" placeholder in B3.
# Only the placeholder should survive (and it should lose its border too);
# the other sample cells (C3:E3) are removed outright.
$ws.Range("C3:E3").Clear()
$ws.Range("B3").Borders.LineStyle = 0

# Row 21 held a second placeholder ("[Table2]") that's no longer needed.
$ws.Range("B21").Clear()

# Leave the selection sitting on the remaining placeholder cell.
$null = $ws.Range("B3").Select()
